$wb = $excel.ActiveWorkbook

# Sheet 1: Overview
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("G4").Value = "2016-08-30 06:17:56"
$wsOverview.Range("G5").Value = "2016-08-30 06:17:56"

# Sheet 2: zh-cn
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-30 06:17:51"
$wsZhCn.Range("H5").Value = "2016-08-30 06:17:51"
$wsZhCn.Range("K4").Value = "2016-08-30 06:18:15"
$wsZhCn.Range("K5").Value = "2016-08-30 06:18:15"

# Sheet 3: de-de
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-30 06:17:56"
$wsDeDe.Range("H5").Value = "2016-08-30 06:17:56"
$wsDeDe.Range("K4").Value = "2016-08-30 06:18:22"
$wsDeDe.Range("K5").Value = "2016-08-30 06:18:22"
